$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("W")

# Row 6: Change in inventories
$ws.Range("B6").Value = 9947000.0
$ws.Range("C6").Value = 14687000.0
$ws.Range("D6").Value = 2657000.0
$ws.Range("E6").Value = -13035000.0
$ws.Range("F6").Value = -15631000.0

# Row 7: Change in payables and accrued liability
$ws.Range("B7").Value = 531526000.0
$ws.Range("C7").Value = 1264215000.0
$ws.Range("D7").Value = 1364841000.0
$ws.Range("E7").Value = 962532000.0
$ws.Range("F7").Value = 1123786000.0
